$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 148, shifting existing data
# (old rows 148-166) down to rows 150-168.
$ws.Rows.Item(148).Insert()
$ws.Rows.Item(148).Insert()

# New row 148: Valencia / Primera entry dated 2021-11-05 (serial 44505)
$ws.Cells.Item(148,1).Value  = 11
$ws.Cells.Item(148,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(148,3).Value  = "Bíobío"
$ws.Cells.Item(148,4).Value  = 44505
$ws.Cells.Item(148,5).Value  = 8
$ws.Cells.Item(148,6).Value  = "Fruta"
$ws.Cells.Item(148,7).Value  = 100102
$ws.Cells.Item(148,8).Value  = "Cítricos"
$ws.Cells.Item(148,9).Value  = 100102005
$ws.Cells.Item(148,10).Value = "Naranja"
$ws.Cells.Item(148,11).Value = "Valencia"
$ws.Cells.Item(148,12).Value = "Primera"
$ws.Cells.Item(148,13).Value = 200
$ws.Cells.Item(148,14).Value = 8000
$ws.Cells.Item(148,15).Value = 9000
$ws.Cells.Item(148,16).Value = 8500
$ws.Cells.Item(148,17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(148,18).Value = "Región de O'Higgins"
$ws.Cells.Item(148,19).Value = 567
$ws.Cells.Item(148,20).Value = 15

# New row 149: Valencia / Segunda entry, same date
$ws.Cells.Item(149,1).Value  = 11
$ws.Cells.Item(149,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(149,3).Value  = "Bíobío"
$ws.Cells.Item(149,4).Value  = 44505
$ws.Cells.Item(149,5).Value  = 8
$ws.Cells.Item(149,6).Value  = "Fruta"
$ws.Cells.Item(149,7).Value  = 100102
$ws.Cells.Item(149,8).Value  = "Cítricos"
$ws.Cells.Item(149,9).Value  = 100102005
$ws.Cells.Item(149,10).Value = "Naranja"
$ws.Cells.Item(149,11).Value = "Valencia"
$ws.Cells.Item(149,12).Value = "Segunda"
$ws.Cells.Item(149,13).Value = 100
$ws.Cells.Item(149,14).Value = 7000
$ws.Cells.Item(149,15).Value = 7000
$ws.Cells.Item(149,16).Value = 7000
$ws.Cells.Item(149,17).Value = "`$/caja 15 kilos empedrada"
$ws.Cells.Item(149,18).Value = "Región de O'Higgins"
$ws.Cells.Item(149,19).Value = 467
$ws.Cells.Item(149,20).Value = 15
